$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set cell values
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Format B1: bold font, thin border all around, centered horizontally, top-aligned vertically
$cell = $ws.Range("B1")
$cell.HorizontalAlignment = -4108  # xlCenter
$cell.VerticalAlignment = -4160    # xlTop
$cell.Font.Bold = $true
$cell.Borders.LineStyle = 1        # xlContinuous
$cell.Borders.Weight = 2           # xlThin

# Copy the same formatting onto A2
$ws.Range("B1").Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
